$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Cédula column values (A2:A8) with malformed/odd formats that the
#     new validation rule is meant to catch ---
$ws.Range("A2").Value = "V-265541DJNDJNJ"
$ws.Range("A3").Value = "EE26554126"
$ws.Range("A4").Value = "V-26554127"

# A5 becomes a bare number (left aligned) instead of a "V-xxxxx" text id
$ws.Range("A5").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("A5").Value = 5662332

$ws.Range("A6").Value = "V26554129"
$ws.Range("A7").Value = "v26554130"
$ws.Range("A8").Value = "V- 26554131"

# --- Carrera column fix for row 7 ---
$ws.Range("D7").Value = "Biología"

# --- New student rows 9 and 10 ---
$ws.Range("B9").Value = "Estudiante 8"
$ws.Range("C9").Value = "Apellido 8"
$ws.Range("D9").Value = "Computación"

# Row 10: A10 is literal text consisting of apostrophes. Typing a leading
# apostrophe in Excel marks the cell as quote-prefixed (text literal) and is
# not stored as part of the value, so 8 apostrophes here yields the 7
# stored apostrophes + quotePrefix style seen in the target file.
$ws.Range("A10").Value = "''''''''"
$ws.Range("B10").Value = "Estudiante 9"
$ws.Range("C10").Value = "Apellido 9"
$ws.Range("D10").Value = "Computación"

$ws.Range("E10").Select()

$wb.Save()
